$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.414.03"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.608.42"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'212.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.0607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'19.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.836.66"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "1.605.96"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "'63.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "'233.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.27%  "
$ws.Range("D18").Value = "26.419.05"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'7.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "'8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "'147.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'6.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "1.491.78"
$ws.Range("E32").Value = "  +5.23%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").Value = "'0.939"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").Value = "1.747.44"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'60.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'89.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").Value = "'1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'0.0501"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'0.0964"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("E51").Value = "  +1.54%  "
